## Agrega las especificaciones del Sprint #5 (y Sprint 4, y la seccion
## "Desafios del proyecto") al documento de especificaciones del proyecto.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Bloque 1: Sprint 4 y Sprint 5, insertado justo despues del primer parrafo
# vacio que sigue a "Crear la ventana principal de la Interfaz Grafica."
# (es decir, entre los dos parrafos vacios que preceden a "Historias de
# usuario").
# ---------------------------------------------------------------------------

$sprintsXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Sprint 4 – Crear las funcionalidades para guardar libros (IGU, lógica y persistencia) (1</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Semana) 22 – 28 enero 2024</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Crear la interfaz gráfica para guardar libros (capa presentación).</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Crear la lógica para crear libros (capa lógica).</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Crear la</w:t></w:r><w:r><w:t xml:space="preserve"> sentencia SQL para guardar libros</w:t></w:r><w:r><w:t xml:space="preserve"> en la base de datos</w:t></w:r><w:r><w:t xml:space="preserve"> (capa persistencia). </w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Sprint 5 – Traer los datos de los libros guardados (SELECT) y mostrarlos en una tabla (1 Semana) 29 enero – 04 febrero 2024</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Crear una tabla para mostrar los datos de los libros guardados.</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Crear la lógica que permita solicitar los datos a la persistencia y convertir al tipo de dato adecuado.</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Crear la sentencia SQL para solicitar los libros almacenados en la tabla libros. </w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@

$anchor1 = $d.Content
$anchor1.Find.Execute("Crear la ventana principal de la Interfaz Gráfica.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorPara1 = $anchor1.Paragraphs(1)
$firstEmptyPara = $anchorPara1.Next()
$insertRange1 = $firstEmptyPara.Range
$insertRange1.Collapse(0)
$insertRange1.InsertXML($sprintsXml)

# ---------------------------------------------------------------------------
# Bloque 2: seccion "Desafios del proyecto", insertado justo despues del
# parrafo vacio que sigue a "Historias de usuario".
# ---------------------------------------------------------------------------

$desafiosXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Ttulo2"/></w:pPr><w:r><w:t>Desafíos del proyecto</w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">En esta sección se incorporan los desafíos con los que he tenido que lidiar al momento de realizar este proyecto. </w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Conversión de una imagen a tipo de dato </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>byte[</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">] (arreglo de bytes) para almacenar en la base de datos. </w:t></w:r></w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@

$anchor2 = $d.Content
$anchor2.Find.Execute("Historias de usuario", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorPara2 = $anchor2.Paragraphs(1)
$secondEmptyPara = $anchorPara2.Next()
$insertRange2 = $secondEmptyPara.Range
$insertRange2.Collapse(0)
$insertRange2.InsertXML($desafiosXml)

Write-Output "done"
